$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1409.3462
$ws.Range("I19").Value = 845.25
$ws.Range("J19").Value = 1892.8572
$ws.Range("K19").Value = 845.25
$ws.Range("L19").Value = 1892.8572
$ws.Range("M19").Value = -670.25
$ws.Range("N19").Value = -2242.8572
$ws.Range("H48").Value = 1988
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1988
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 5964
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -6548
$ws.Range("H56").Value = 1988
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1988
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 5964
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -7032
$ws.Range("H74").Value = 6388.6
$ws.Range("I74").Value = 6388.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6388.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5452.6
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 6388.6
$ws.Range("I77").Value = 6388.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 31943
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -27263
$ws.Range("N77").ClearContents()
$ws.Range("H112").Value = 3032223
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 3670270
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 11010810
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -11013026
$ws.Range("H129").Value = 1753
$ws.Range("I129").Value = 1753
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 5259
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -259
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 32222.975
$ws.Range("I137").Value = 67285.766
$ws.Range("J137").Value = 3838.8096
$ws.Range("K137").Value = 201857.298
$ws.Range("L137").Value = 11516.4288
$ws.Range("M137").Value = -199307.298
$ws.Range("N137").Value = -16616.4288
$ws.Range("H138").Value = 1879.3334
$ws.Range("I138").Value = 1879.3334
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 5638.0002
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -498.0002000000004
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 2046.0435
$ws.Range("I141").Value = 2120.9048
$ws.Range("K141").Value = 6362.714399999999
$ws.Range("M141").Value = -1182.714399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2633.353
$ws.Range("I61").Value = 2502.2068
$ws.Range("J61").Value = 3394
$ws.Range("K61").Value = 2502.2068
$ws.Range("L61").Value = 3394
$ws.Range("M61").Value = -2290.2068
$ws.Range("N61").Value = -3818
$ws.Range("H63").Value = 3639.4
$ws.Range("I63").Value = 3639.4
$ws.Range("K63").Value = 3639.4
$ws.Range("M63").Value = -2953.4
$ws.Range("H66").Value = 3639.4
$ws.Range("I66").Value = 3639.4
$ws.Range("K66").Value = 18197
$ws.Range("M66").Value = -14765
$ws.Range("H122").Value = 2634.5806
$ws.Range("I122").Value = 2609.8276
$ws.Range("J122").Value = 2993.5
$ws.Range("K122").Value = 7829.4828
$ws.Range("L122").Value = 8980.5
$ws.Range("M122").Value = -5379.4828
$ws.Range("N122").Value = -13880.5
$ws.Range("H132").Value = 3115.6667
$ws.Range("I132").Value = 3115.6667
$ws.Range("K132").Value = 9347.000100000001
$ws.Range("M132").Value = -6817.000100000001
$ws.Range("H136").Value = 2633.353
$ws.Range("I136").Value = 2502.2068
$ws.Range("J136").Value = 3394
$ws.Range("K136").Value = 7506.6204
$ws.Range("L136").Value = 10182
$ws.Range("M136").Value = -4956.6204
$ws.Range("N136").Value = -15282

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1818.5883
$ws.Range("I105").Value = 1804.7858
$ws.Range("J105").Value = 1883
$ws.Range("K105").Value = 1804.7858
$ws.Range("L105").Value = 1883
$ws.Range("M105").Value = -57.78580000000011
$ws.Range("N105").Value = -5377
$ws.Range("H134").Value = 4210.911
$ws.Range("I134").Value = 5154.8623
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 15464.5869
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -12929.5869
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 145569.44
$ws.Range("I31").Value = 214940.16
$ws.Range("J31").Value = 3811.913
$ws.Range("K31").Value = 214940.16
$ws.Range("L31").Value = 3811.913
$ws.Range("M31").Value = -214645.16
$ws.Range("N31").Value = -4401.913
$ws.Range("H34").Value = 145569.44
$ws.Range("I34").Value = 214940.16
$ws.Range("J34").Value = 3811.913
$ws.Range("K34").Value = 214940.16
$ws.Range("L34").Value = 3811.913
$ws.Range("M34").Value = -214738.16
$ws.Range("N34").Value = -4215.913
$ws.Range("H58").Value = 4045.1333
$ws.Range("I58").Value = 3133
$ws.Range("J58").Value = 5869.4
$ws.Range("K58").Value = 3133
$ws.Range("L58").Value = 5869.4
$ws.Range("M58").Value = -2930
$ws.Range("N58").Value = -6275.4
$ws.Range("H107").Value = 3950.5715
$ws.Range("I107").Value = 786.1579
$ws.Range("J107").Value = 6564.6523
$ws.Range("K107").Value = 786.1579
$ws.Range("L107").Value = 6564.6523
$ws.Range("M107").Value = 1133.8421
$ws.Range("N107").Value = -10404.6523
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 4767.893
$ws.Range("I132").Value = 4279.933
$ws.Range("J132").Value = 6764.091
$ws.Range("K132").Value = 12839.799
$ws.Range("L132").Value = 20292.273
$ws.Range("M132").Value = -10309.799
$ws.Range("N132").Value = -25352.273
$ws.Range("H134").Value = 2291.946
$ws.Range("I134").Value = 2359.457
$ws.Range("K134").Value = 7078.370999999999
$ws.Range("M134").Value = -4543.370999999999
$ws.Range("H136").Value = 4045.1333
$ws.Range("I136").Value = 3133
$ws.Range("J136").Value = 5869.4
$ws.Range("K136").Value = 9399
$ws.Range("L136").Value = 17608.2
$ws.Range("M136").Value = -6849
$ws.Range("N136").Value = -22708.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3481.4285
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 3995
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 11985
$ws.Range("M32").Value = -917
$ws.Range("N32").Value = -12551
$ws.Range("H46").Value = 1167.091
$ws.Range("I46").Value = 991.8571
$ws.Range("J46").Value = 1473.75
$ws.Range("K46").Value = 2975.5713
$ws.Range("L46").Value = 4421.25
$ws.Range("M46").Value = -2884.5713
$ws.Range("N46").Value = -4603.25
$ws.Range("H60").Value = 250
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H121").Value = 111225660
$ws.Range("I121").Value = 200002100
$ws.Range("J121").Value = 255123.75
$ws.Range("K121").Value = 600006300
$ws.Range("L121").Value = 765371.25
$ws.Range("M121").Value = -600004990
$ws.Range("N121").Value = -767991.25
$ws.Range("H131").Value = 50000970
$ws.Range("I131").Value = 50000970
$ws.Range("K131").Value = 150002910
$ws.Range("M131").Value = -149997870
$ws.Range("H132").Value = 6495.619
$ws.Range("I132").Value = 8477.267
$ws.Range("K132").Value = 76295.40299999999
$ws.Range("M132").Value = -73765.40299999999
$ws.Range("H137").Value = 1665.4
$ws.Range("I137").Value = 1248.5
$ws.Range("J137").Value = 3333
$ws.Range("K137").Value = 3745.5
$ws.Range("L137").Value = 9999
$ws.Range("M137").Value = 1354.5
$ws.Range("N137").Value = -20199

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1575.4762
$ws.Range("I113").Value = 1421.3889
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1421.3889
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 748.6111000000001
$ws.Range("N113").Value = -6840
$ws.Range("H122").Value = 4405.933
$ws.Range("I122").Value = 4953.5884
$ws.Range("J122").Value = 3689.7693
$ws.Range("K122").Value = 14860.7652
$ws.Range("L122").Value = 11069.3079
$ws.Range("M122").Value = -12410.7652
$ws.Range("N122").Value = -15969.3079
$ws.Range("H132").Value = 28679.488
$ws.Range("I132").Value = 32220.611
$ws.Range("J132").Value = 3183.4
$ws.Range("K132").Value = 96661.833
$ws.Range("L132").Value = 9550.2
$ws.Range("M132").Value = -94131.833
$ws.Range("N132").Value = -14610.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3322.9583
$ws.Range("I132").Value = 3387.1
$ws.Range("J132").Value = 3002.25
$ws.Range("K132").Value = 10161.3
$ws.Range("L132").Value = 9006.75
$ws.Range("M132").Value = -7631.299999999999
$ws.Range("N132").Value = -14066.75
$ws.Range("H136").Value = 3776.3928
$ws.Range("I136").Value = 3261.56
$ws.Range("K136").Value = 9784.68
$ws.Range("M136").Value = -7234.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("M18").Value = -9827
$ws.Range("H64").Value = 46249.75
$ws.Range("H67").Value = 46249.75
$ws.Range("H112").Value = 60999.668
$ws.Range("J112").Value = 60999.668
$ws.Range("L112").Value = 60999.668
$ws.Range("N112").Value = -63953.668
$ws.Range("H132").Value = 2143.5557
$ws.Range("I132").Value = 1146.2
$ws.Range("J132").Value = 2642.2334
$ws.Range("K132").Value = 3438.6
$ws.Range("L132").Value = 7926.7002
$ws.Range("M132").Value = -908.6000000000004
$ws.Range("N132").Value = -12986.7002
$ws.Range("H136").Value = 358283.2
$ws.Range("I136").Value = 385612.84
$ws.Range("J136").Value = 2997.5
$ws.Range("K136").Value = 1156838.52
$ws.Range("L136").Value = 8992.5
$ws.Range("M136").Value = -1154288.52
$ws.Range("N136").Value = -14092.5
